# Daily attendance processing - 2026-01-03 22:31:23
#
# For every data row in the "Recorded By" column (G), when the cell holds a
# comma-separated list of recorder names/emails, swap the first two entries
# (leaving any further entries, e.g. a trailing lowercase "system", in place).
# Cells holding a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Length -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $newVal = [string]::Join(", ", $parts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
